$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2695.125
$ws.Cells.Item(62, 9).Value = 2624.7693
$ws.Cells.Item(62, 10).Value = 3000
$ws.Cells.Item(62, 11).Value = 2624.7693
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 13).Value = -2000.7693
$ws.Cells.Item(62, 14).Value = -4248
$ws.Cells.Item(65, 8).Value = 2695.125
$ws.Cells.Item(65, 9).Value = 2624.7693
$ws.Cells.Item(65, 10).Value = 3000
$ws.Cells.Item(65, 11).Value = 13123.8465
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = -10003.8465
$ws.Cells.Item(65, 14).Value = -21240
$ws.Cells.Item(107, 8).Value = 1014.9286
$ws.Cells.Item(107, 9).Value = 1023.7692
$ws.Cells.Item(107, 10).Value = 900
$ws.Cells.Item(107, 11).Value = 1023.7692
$ws.Cells.Item(107, 12).Value = 900
$ws.Cells.Item(107, 13).Value = 896.2308
$ws.Cells.Item(107, 14).Value = -4740
$ws.Cells.Item(116, 8).Value = 2212.8572
$ws.Cells.Item(116, 9).Value = 2496.6667
$ws.Cells.Item(116, 11).Value = 2496.6667
$ws.Cells.Item(116, 13).Value = 945.3332999999998
$ws.Cells.Item(129, 8).Value = 2494.951
$ws.Cells.Item(129, 10).Value = 908.4528
$ws.Cells.Item(129, 12).Value = 2725.3584
$ws.Cells.Item(129, 14).Value = -12725.3584
$ws.Cells.Item(132, 8).Value = 4470230
$ws.Cells.Item(132, 9).Value = 5005988
$ws.Cells.Item(132, 10).Value = 5583.1665
$ws.Cells.Item(132, 11).Value = 15017964
$ws.Cells.Item(132, 12).Value = 16749.4995
$ws.Cells.Item(132, 13).Value = -15015434
$ws.Cells.Item(132, 14).Value = -21809.4995
$ws.Cells.Item(137, 8).Value = 1437.2142
$ws.Cells.Item(137, 9).Value = 1240.0769
$ws.Cells.Item(137, 10).Value = 4000
$ws.Cells.Item(137, 11).Value = 3720.2307
$ws.Cells.Item(137, 12).Value = 12000
$ws.Cells.Item(137, 13).Value = -1170.2307
$ws.Cells.Item(137, 14).Value = -17100
$ws.Cells.Item(138, 8).Value = 2898.5176
$ws.Cells.Item(138, 9).Value = 1687.8214
$ws.Cells.Item(138, 10).Value = 3493.2456
$ws.Cells.Item(138, 11).Value = 5063.4642
$ws.Cells.Item(138, 12).Value = 10479.7368
$ws.Cells.Item(138, 13).Value = 76.53579999999965
$ws.Cells.Item(138, 14).Value = -20759.7368
$ws.Cells.Item(141, 8).Value = 3556.25
$ws.Cells.Item(141, 9).Value = 3252.7778
$ws.Cells.Item(141, 10).Value = 4466.6665
$ws.Cells.Item(141, 11).Value = 9758.3334
$ws.Cells.Item(141, 12).Value = 13399.9995
$ws.Cells.Item(141, 13).Value = -4578.3334
$ws.Cells.Item(141, 14).Value = -23759.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 201251.4
$ws.Cells.Item(2, 9).Value = 1564.25
$ws.Cells.Item(2, 10).Value = 1000000
$ws.Cells.Item(2, 11).Value = 1564.25
$ws.Cells.Item(2, 12).Value = 1000000
$ws.Cells.Item(2, 13).Value = -1451.25
$ws.Cells.Item(2, 14).Value = -1000226
$ws.Cells.Item(10, 8).Value = 2000
$ws.Cells.Item(10, 10).Value = 2000
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 14).Value = -2340
$ws.Cells.Item(23, 8).Value = 12000
$ws.Cells.Item(23, 10).Value = 12000
$ws.Cells.Item(23, 12).Value = 12000
$ws.Cells.Item(23, 14).Value = -12518
$ws.Cells.Item(32, 8).Value = 29775.412
$ws.Cells.Item(32, 9).Value = 12004.429
$ws.Cells.Item(32, 10).Value = 80665.05
$ws.Cells.Item(32, 11).Value = 12004.429
$ws.Cells.Item(32, 12).Value = 80665.05
$ws.Cells.Item(32, 13).Value = -11717.429
$ws.Cells.Item(32, 14).Value = -81239.05
$ws.Cells.Item(37, 8).Value = 12250
$ws.Cells.Item(37, 9).Value = 6800
$ws.Cells.Item(37, 10).Value = 14975
$ws.Cells.Item(37, 11).Value = 6800
$ws.Cells.Item(37, 12).Value = 14975
$ws.Cells.Item(37, 13).Value = -6527
$ws.Cells.Item(37, 14).Value = -15521
$ws.Cells.Item(45, 8).Value = 101616.5
$ws.Cells.Item(45, 9).Value = 250952.75
$ws.Cells.Item(45, 10).Value = 2059
$ws.Cells.Item(45, 11).Value = 250952.75
$ws.Cells.Item(45, 12).Value = 2059
$ws.Cells.Item(45, 13).Value = -250575.75
$ws.Cells.Item(45, 14).Value = -2813
$ws.Cells.Item(63, 8).Value = 2567
$ws.Cells.Item(63, 9).Value = 1725.75
$ws.Cells.Item(63, 10).Value = 3240
$ws.Cells.Item(63, 11).Value = 1725.75
$ws.Cells.Item(63, 12).Value = 3240
$ws.Cells.Item(63, 13).Value = -1039.75
$ws.Cells.Item(63, 14).Value = -4612
$ws.Cells.Item(66, 8).Value = 2567
$ws.Cells.Item(66, 9).Value = 1725.75
$ws.Cells.Item(66, 10).Value = 3240
$ws.Cells.Item(66, 11).Value = 8628.75
$ws.Cells.Item(66, 12).Value = 16200
$ws.Cells.Item(66, 13).Value = -5196.75
$ws.Cells.Item(66, 14).Value = -23064
$ws.Cells.Item(74, 8).Value = 1564.1904
$ws.Cells.Item(74, 9).Value = 1652.6428
$ws.Cells.Item(74, 10).Value = 1387.2858
$ws.Cells.Item(74, 11).Value = 1652.6428
$ws.Cells.Item(74, 12).Value = 1387.2858
$ws.Cells.Item(74, 13).Value = -778.6428000000001
$ws.Cells.Item(74, 14).Value = -3135.2858
$ws.Cells.Item(77, 8).Value = 1564.1904
$ws.Cells.Item(77, 9).Value = 1652.6428
$ws.Cells.Item(77, 10).Value = 1387.2858
$ws.Cells.Item(77, 11).Value = 8263.214
$ws.Cells.Item(77, 12).Value = 6936.429
$ws.Cells.Item(77, 13).Value = -3895.214
$ws.Cells.Item(77, 14).Value = -15672.429
$ws.Cells.Item(116, 8).Value = 201251.4
$ws.Cells.Item(116, 9).Value = 1564.25
$ws.Cells.Item(116, 10).Value = 1000000
$ws.Cells.Item(116, 11).Value = 1564.25
$ws.Cells.Item(116, 12).Value = 1000000
$ws.Cells.Item(116, 13).Value = 729.75
$ws.Cells.Item(116, 14).Value = -1004588

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 201251.4
$ws.Cells.Item(3, 9).Value = 1564.25
$ws.Cells.Item(3, 10).Value = 1000000
$ws.Cells.Item(3, 11).Value = 1564.25
$ws.Cells.Item(3, 12).Value = 1000000
$ws.Cells.Item(3, 13).Value = -1450.25
$ws.Cells.Item(3, 14).Value = -1000228
$ws.Cells.Item(7, 8).Value = 451.5
$ws.Cells.Item(7, 9).Value = 451.5
$ws.Cells.Item(7, 11).Value = 451.5
$ws.Cells.Item(7, 13).Value = -338.5
$ws.Cells.Item(20, 8).Value = 48845.727
$ws.Cells.Item(20, 9).Value = 59272.555
$ws.Cells.Item(20, 10).Value = 1925
$ws.Cells.Item(20, 11).Value = 59272.555
$ws.Cells.Item(20, 12).Value = 1925
$ws.Cells.Item(20, 13).Value = -59025.555
$ws.Cells.Item(20, 14).Value = -2419
$ws.Cells.Item(35, 8).Value = 20797.428
$ws.Cells.Item(35, 10).Value = 21430.334
$ws.Cells.Item(35, 12).Value = 21430.334
$ws.Cells.Item(35, 14).Value = -22050.334
$ws.Cells.Item(141, 8).Value = 53675

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 42568.64
$ws.Cells.Item(31, 9).Value = 1296.4375
$ws.Cells.Item(31, 10).Value = 75586.39999999999
$ws.Cells.Item(31, 11).Value = 1296.4375
$ws.Cells.Item(31, 12).Value = 75586.39999999999
$ws.Cells.Item(31, 13).Value = -1001.4375
$ws.Cells.Item(31, 14).Value = -76176.39999999999
$ws.Cells.Item(34, 8).Value = 42568.64
$ws.Cells.Item(34, 9).Value = 1296.4375
$ws.Cells.Item(34, 10).Value = 75586.39999999999
$ws.Cells.Item(34, 11).Value = 1296.4375
$ws.Cells.Item(34, 12).Value = 75586.39999999999
$ws.Cells.Item(34, 13).Value = -1094.4375
$ws.Cells.Item(34, 14).Value = -75990.39999999999
$ws.Cells.Item(58, 8).Value = 1504.35
$ws.Cells.Item(58, 9).Value = 1284.8387
$ws.Cells.Item(58, 10).Value = 2260.4443
$ws.Cells.Item(58, 11).Value = 1284.8387
$ws.Cells.Item(58, 12).Value = 2260.4443
$ws.Cells.Item(58, 13).Value = -1081.8387
$ws.Cells.Item(58, 14).Value = -2666.4443
$ws.Cells.Item(94, 8).Value = 1284.7273
$ws.Cells.Item(94, 9).Value = 1218.5
$ws.Cells.Item(94, 10).Value = 1309.5625
$ws.Cells.Item(94, 11).Value = 1218.5
$ws.Cells.Item(94, 12).Value = 1309.5625
$ws.Cells.Item(94, 13).Value = -767.5
$ws.Cells.Item(94, 14).Value = -2211.5625
$ws.Cells.Item(136, 8).Value = 1504.35
$ws.Cells.Item(136, 9).Value = 1284.8387
$ws.Cells.Item(136, 10).Value = 2260.4443
$ws.Cells.Item(136, 11).Value = 3854.5161
$ws.Cells.Item(136, 12).Value = 6781.3329
$ws.Cells.Item(136, 13).Value = -1304.5161
$ws.Cells.Item(136, 14).Value = -11881.3329

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 1837
$ws.Cells.Item(63, 9).Value = 1837
$ws.Cells.Item(63, 11).Value = 5511
$ws.Cells.Item(63, 13).Value = -4762
$ws.Cells.Item(66, 8).Value = 1837
$ws.Cells.Item(66, 9).Value = 1837
$ws.Cells.Item(66, 11).Value = 16533
$ws.Cells.Item(66, 13).Value = -12789
$ws.Cells.Item(121, 8).Value = 4193.129
$ws.Cells.Item(121, 9).Value = 11747.25
$ws.Cells.Item(121, 10).Value = 3074
$ws.Cells.Item(121, 11).Value = 35241.75
$ws.Cells.Item(121, 12).Value = 9222
$ws.Cells.Item(121, 13).Value = -33931.75
$ws.Cells.Item(121, 14).Value = -11842
$ws.Cells.Item(132, 8).Value = 456210.12
$ws.Cells.Item(132, 9).Value = 810.36365
$ws.Cells.Item(132, 10).Value = 911609.9399999999
$ws.Cells.Item(132, 11).Value = 7293.27285
$ws.Cells.Item(132, 12).Value = 8204489.459999999
$ws.Cells.Item(132, 13).Value = -4763.27285
$ws.Cells.Item(132, 14).Value = -8209549.459999999
$ws.Cells.Item(140, 8).Value = 4837.8276
$ws.Cells.Item(140, 9).Value = 5251.88
$ws.Cells.Item(140, 11).Value = 15755.64
$ws.Cells.Item(140, 13).Value = -10575.64

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 4236.222
$ws.Cells.Item(55, 10).Value = 5344.3335
$ws.Cells.Item(55, 12).Value = 5344.3335
$ws.Cells.Item(55, 14).Value = -5998.3335
$ws.Cells.Item(122, 8).Value = 1553.4615
$ws.Cells.Item(122, 9).Value = 1770
$ws.Cells.Item(122, 10).Value = 831.6667
$ws.Cells.Item(122, 11).Value = 5310
$ws.Cells.Item(122, 12).Value = 2495.0001
$ws.Cells.Item(122, 13).Value = -2860
$ws.Cells.Item(122, 14).Value = -7395.0001
$ws.Cells.Item(126, 8).Value = 5885423
$ws.Cells.Item(126, 9).Value = 3200.5
$ws.Cells.Item(126, 10).Value = 14708757
$ws.Cells.Item(126, 11).Value = 9601.5
$ws.Cells.Item(126, 12).Value = 44126271
$ws.Cells.Item(126, 13).Value = -7131.5
$ws.Cells.Item(126, 14).Value = -44131211
$ws.Cells.Item(135, 8).Value = 44985
$ws.Cells.Item(135, 10).Value = 44985
$ws.Cells.Item(135, 12).Value = 44985
$ws.Cells.Item(135, 14).Value = -55125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2396.6667
$ws.Cells.Item(7, 9).Value = 1768.75
$ws.Cells.Item(7, 10).Value = 3114.2856
$ws.Cells.Item(7, 11).Value = 1768.75
$ws.Cells.Item(7, 12).Value = 3114.2856
$ws.Cells.Item(7, 13).Value = -1656.75
$ws.Cells.Item(7, 14).Value = -3338.2856
$ws.Cells.Item(40, 8).Value = 74254.07000000001
$ws.Cells.Item(40, 9).Value = 169001.17
$ws.Cells.Item(40, 10).Value = 3193.75
$ws.Cells.Item(40, 11).Value = 169001.17
$ws.Cells.Item(40, 12).Value = 3193.75
$ws.Cells.Item(40, 13).Value = -168865.17
$ws.Cells.Item(40, 14).Value = -3465.75
$ws.Cells.Item(61, 8).Value = 4099.857
$ws.Cells.Item(61, 9).Value = 3100.889
$ws.Cells.Item(61, 10).Value = 5898
$ws.Cells.Item(61, 11).Value = 3100.889
$ws.Cells.Item(61, 12).Value = 5898
$ws.Cells.Item(61, 13).Value = -2898.889
$ws.Cells.Item(61, 14).Value = -6302
$ws.Cells.Item(113, 8).Value = 4099.857
$ws.Cells.Item(113, 9).Value = 3100.889
$ws.Cells.Item(113, 10).Value = 5898
$ws.Cells.Item(113, 11).Value = 3100.889
$ws.Cells.Item(113, 12).Value = 5898
$ws.Cells.Item(113, 13).Value = -930.8890000000001
$ws.Cells.Item(113, 14).Value = -10238
$ws.Cells.Item(126, 8).Value = 2396.6667
$ws.Cells.Item(126, 9).Value = 1768.75
$ws.Cells.Item(126, 10).Value = 3114.2856
$ws.Cells.Item(126, 11).Value = 5306.25
$ws.Cells.Item(126, 12).Value = 9342.856800000001
$ws.Cells.Item(126, 13).Value = -2836.25
$ws.Cells.Item(126, 14).Value = -14282.8568
$ws.Cells.Item(136, 8).Value = 2318.375
$ws.Cells.Item(136, 9).Value = 1522.1305
$ws.Cells.Item(136, 10).Value = 4353.222
$ws.Cells.Item(136, 11).Value = 4566.3915
$ws.Cells.Item(136, 12).Value = 13059.666
$ws.Cells.Item(136, 13).Value = -2016.3915
$ws.Cells.Item(136, 14).Value = -18159.666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 18555.555
$ws.Cells.Item(123, 10).Value = 18555.555
$ws.Cells.Item(123, 12).Value = 18555.555
$ws.Cells.Item(123, 14).Value = -28355.555
$ws.Cells.Item(136, 8).Value = 24878.5
$ws.Cells.Item(136, 9).Value = 72140.57000000001
$ws.Cells.Item(136, 10).Value = 5417.647
$ws.Cells.Item(136, 11).Value = 216421.71
$ws.Cells.Item(136, 12).Value = 16252.941
$ws.Cells.Item(136, 13).Value = -213871.71
$ws.Cells.Item(136, 14).Value = -21352.941
